$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Challenges that have now been passed (rows 123-128): update their text from
# "Not Passed..." to "Passed..." and record the completion date.
$ws.Range("D123").Value = "PassedUse Destructuring Assignment to Assign Variables from Nested Objects"
$ws.Range("D124").Value = "PassedUse Destructuring Assignment to Assign Variables from Arrays"
$ws.Range("D125").Value = "PassedUse Destructuring Assignment with the Rest Operator to Reassign Array Elements"
$ws.Range("D126").Value = "PassedUse Destructuring Assignment to Pass an Object as a Function's Parameters"
$ws.Range("D127").Value = "PassedCreate Strings using Template Literals"
$ws.Range("D128").Value = "PassedWrite Concise Object Literal Declarations Using Simple Fields"

# Rows 129 and 130 are also completed now (text/title stays the same), but they
# still need their completion date recorded and highlight cleared below.
$dateDone = 43437
$ws.Range("E123:E130").Value = $dateDone

# Clear the "currently working on" highlight (bold + fill) from the rows that
# are now finished.
$doneRows = $ws.Range("D123:D130")
$doneRows.Font.Bold = $false
$doneRows.Interior.Pattern = -4142

# Move the "currently working on" highlight down to the next two challenges.
$nextUp = $ws.Range("D131:D132")
$nextUp.Interior.Color = 65535

# Update the active selection to reflect where editing left off.
$ws.Range("E130").Select()
